# "Generate Report for handoff"
#
# Regenerate the localization-status report:
#  - the previously "Handoff transform failed" entry (GUID
#    6286a789-0224-4352-b49a-cd3226406aca) is gone, its row is removed and
#    the ".localization-config" row shifts up to take its place
#  - a new handoff (GUID 9a5f9c62-924d-4520-99cb-4fba81d62124, content hash
#    db5437abdaf804858dd245424949c9d53215c572) replaces the old
#    bb85d3f6-00ee-4245-8fdf-c04ac291d568 entry, with fresh handoff
#    timestamps for zh-cn and de-de

$wb = $excel.ActiveWorkbook

$oldGuid = "bb85d3f6-00ee-4245-8fdf-c04ac291d568"
$newGuid = "9a5f9c62-924d-4520-99cb-4fba81d62124"
$newHash = "db5437abdaf804858dd245424949c9d53215c572"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newZhTimestamp = "2016-01-18 04:04:24"
$newDeTimestamp = "2016-01-18 04:04:35"

$repoBase    = "https://github.com/OpenLocalizationTest/oltest/blob/676acdc4b416cd6ca1ce18ee9753475ea935c127"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/133b9a8a14bd6d9914b52ac8cb6f3ed5c08c5400/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b88327498073f437143b779b59dcb610d0e8203b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"

$newMdUrl    = "$repoBase/e2e/$newMdName"
$configUrl   = "$repoBase/.localization-config"
$newZhXlfUrl = "$zhHandoffBase/$newZhXlfName"
$newDeXlfUrl = "$deHandoffBase/$newDeXlfName"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Drop the old "Handoff transform failed" row (row 3); the
# ".localization-config" row below it shifts up into row 3.
$ws.Rows.Item(3).Delete()

# Row 2 now describes the new handoff.
$ws.Range("A2").Value2 = $newMdName

# Rebuild the hyperlinks for this sheet against the new row layout.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value2 = $newMdName
$ws.Range("C2").Value2 = $newZhXlfName
$ws.Range("D2").Value2 = $newZhTimestamp

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $newZhXlfUrl, [Type]::Missing, [Type]::Missing, $newZhXlfName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value2 = $newMdName
$ws.Range("C2").Value2 = $newDeXlfName
$ws.Range("D2").Value2 = $newDeTimestamp

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newMdUrl, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $newDeXlfUrl, [Type]::Missing, [Type]::Missing, $newDeXlfName) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
